# Refresh the market-price-derived columns (H:N — currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) on each job
# sheet with the latest pulled values, as produced by the scheduled
# market-data runner. These are plain cached values (no formulas), so
# each touched cell is written directly.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 183.61539
$ws.Range("I33").Value = 221.11111
$ws.Range("K33").Value = 221.11111
$ws.Range("M33").Value = 7.888890000000004
$ws.Range("H41").Value = 2769.8572
$ws.Range("J41").Value = 3477.6
$ws.Range("L41").Value = 3477.6
$ws.Range("N41").Value = -4357.6
$ws.Range("H70").Value = 2180
$ws.Range("I70").Value = 405
$ws.Range("J70").Value = 3600
$ws.Range("K70").Value = 1215
$ws.Range("L70").Value = 10800
$ws.Range("M70").Value = -945
$ws.Range("N70").Value = -11340
$ws.Range("H73").Value = 2180
$ws.Range("I73").Value = 405
$ws.Range("J73").Value = 3600
$ws.Range("K73").Value = 1215
$ws.Range("L73").Value = 10800
$ws.Range("M73").Value = -279
$ws.Range("N73").Value = -12672
$ws.Range("H76").Value = 3974.625
$ws.Range("I76").Value = 3459.8
$ws.Range("J76").Value = 4832.6665
$ws.Range("K76").Value = 3459.8
$ws.Range("L76").Value = 4832.6665
$ws.Range("M76").Value = -3144.8
$ws.Range("N76").Value = -5462.6665
$ws.Range("H79").Value = 3974.625
$ws.Range("I79").Value = 3459.8
$ws.Range("J79").Value = 4832.6665
$ws.Range("K79").Value = 3459.8
$ws.Range("L79").Value = 4832.6665
$ws.Range("M79").Value = -2367.8
$ws.Range("N79").Value = -7016.6665
$ws.Range("H132").Value = 100009.83
$ws.Range("I132").Value = 203013.16
$ws.Range("K132").Value = 609039.48
$ws.Range("M132").Value = -606509.48
$ws.Range("H137").Value = 2451.5386
$ws.Range("I137").Value = 1597.2759
$ws.Range("J137").Value = 4928.9
$ws.Range("K137").Value = 4791.8277
$ws.Range("L137").Value = 14786.7
$ws.Range("M137").Value = -2241.8277
$ws.Range("N137").Value = -19886.7
$ws.Range("H138").Value = 5411.655
$ws.Range("J138").Value = 7000.8716
$ws.Range("L138").Value = 21002.6148
$ws.Range("N138").Value = -31282.6148
$ws.Range("H140").Value = 68072.22
$ws.Range("J140").Value = 67831.25
$ws.Range("L140").Value = 67831.25
$ws.Range("N140").Value = -78191.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 769
$ws.Range("I5").Value = 769
$ws.Range("K5").Value = 769
$ws.Range("M5").Value = -657
$ws.Range("H32").Value = 3536.623
$ws.Range("I32").Value = 1999.2858
$ws.Range("J32").Value = 9814.083000000001
$ws.Range("K32").Value = 1999.2858
$ws.Range("L32").Value = 9814.083000000001
$ws.Range("M32").Value = -1712.2858
$ws.Range("N32").Value = -10388.083
$ws.Range("H61").Value = 18733.04
$ws.Range("I61").Value = 12207.066
$ws.Range("K61").Value = 12207.066
$ws.Range("M61").Value = -11995.066
$ws.Range("H88").Value = 7574.875
$ws.Range("J88").Value = 9600
$ws.Range("L88").Value = 9600
$ws.Range("N88").Value = -10412
$ws.Range("H91").Value = 7574.875
$ws.Range("J91").Value = 9600
$ws.Range("L91").Value = 9600
$ws.Range("N91").Value = -12408
$ws.Range("H132").Value = 28548.64
$ws.Range("I132").Value = 31327.475
$ws.Range("K132").Value = 93982.42499999999
$ws.Range("M132").Value = -91452.42499999999
$ws.Range("H136").Value = 18733.04
$ws.Range("I136").Value = 12207.066
$ws.Range("K136").Value = 36621.198
$ws.Range("M136").Value = -34071.198

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 769
$ws.Range("I4").Value = 769
$ws.Range("K4").Value = 769
$ws.Range("M4").Value = -654
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H64").Value = 4167564
$ws.Range("I64").Value = 6945249
$ws.Range("K64").Value = 6945249
$ws.Range("M64").Value = -6945024
$ws.Range("H67").Value = 4167564
$ws.Range("I67").Value = 6945249
$ws.Range("K67").Value = 6945249
$ws.Range("M67").Value = -6944469
$ws.Range("H134").Value = 4997.4326
$ws.Range("I134").Value = 2867.7368
$ws.Range("K134").Value = 8603.2104
$ws.Range("M134").Value = -6068.2104
$ws.Range("H140").Value = 225110
$ws.Range("J140").Value = 225110
$ws.Range("L140").Value = 225110
$ws.Range("N140").Value = -235470
$ws.Range("H141").Value = 72624
$ws.Range("J141").Value = 72624
$ws.Range("L141").Value = 72624
$ws.Range("N141").Value = -82984

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 312.22223
$ws.Range("I22").Value = 313.75
$ws.Range("K22").Value = 313.75
$ws.Range("M22").Value = 36.25
$ws.Range("H31").Value = 16397478
$ws.Range("I31").Value = 30305208
$ws.Range("J31").Value = 6223.7856
$ws.Range("K31").Value = 30305208
$ws.Range("L31").Value = 6223.7856
$ws.Range("M31").Value = -30304913
$ws.Range("N31").Value = -6813.7856
$ws.Range("H34").Value = 16397478
$ws.Range("I34").Value = 30305208
$ws.Range("J34").Value = 6223.7856
$ws.Range("K34").Value = 30305208
$ws.Range("L34").Value = 6223.7856
$ws.Range("M34").Value = -30305006
$ws.Range("N34").Value = -6627.7856
$ws.Range("H107").Value = 1015394.3
$ws.Range("I107").Value = 1656304.4
$ws.Range("K107").Value = 1656304.4
$ws.Range("M107").Value = -1654384.4
$ws.Range("H132").Value = 39222484
$ws.Range("I132").Value = 44448050
$ws.Range("K132").Value = 133344150
$ws.Range("M132").Value = -133341620
$ws.Range("H141").Value = 170374.5
$ws.Range("J141").Value = 170374.5
$ws.Range("L141").Value = 170374.5
$ws.Range("N141").Value = -180734.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 289728.6
$ws.Range("J68").Value = 557486.2
$ws.Range("L68").Value = 1672458.6
$ws.Range("N68").Value = -1674080.6
$ws.Range("H71").Value = 289728.6
$ws.Range("J71").Value = 557486.2
$ws.Range("L71").Value = 5017375.8
$ws.Range("N71").Value = -5025487.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 880620.9
$ws.Range("J80").Value = 9281.333000000001
$ws.Range("L80").Value = 9281.333000000001
$ws.Range("N80").Value = -11277.333
$ws.Range("H83").Value = 880620.9
$ws.Range("J83").Value = 9281.333000000001
$ws.Range("L83").Value = 46406.665
$ws.Range("N83").Value = -56390.665
$ws.Range("H102").Value = 6573.75
$ws.Range("I102").Value = 4999.8
$ws.Range("J102").Value = 7698
$ws.Range("K102").Value = 4999.8
$ws.Range("L102").Value = 7698
$ws.Range("M102").Value = -3377.8
$ws.Range("N102").Value = -10942

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4181.4644
$ws.Range("J7").Value = 5615.3076
$ws.Range("L7").Value = 5615.3076
$ws.Range("N7").Value = -5839.3076
$ws.Range("H22").Value = 726.4666999999999
$ws.Range("J22").Value = 730
$ws.Range("L22").Value = 730
$ws.Range("N22").Value = -1320
$ws.Range("H27").Value = 726.4666999999999
$ws.Range("J27").Value = 730
$ws.Range("L27").Value = 730
$ws.Range("N27").Value = -944
$ws.Range("H122").Value = 83344504
$ws.Range("I122").Value = 142861650
$ws.Range("J122").Value = 20511
$ws.Range("K122").Value = 428584950
$ws.Range("L122").Value = 61533
$ws.Range("M122").Value = -428582500
$ws.Range("N122").Value = -66433
$ws.Range("H126").Value = 4181.4644
$ws.Range("J126").Value = 5615.3076
$ws.Range("L126").Value = 16845.9228
$ws.Range("N126").Value = -21785.9228
$ws.Range("H132").Value = 5006.49
$ws.Range("I132").Value = 4300
$ws.Range("J132").Value = 6440.879
$ws.Range("K132").Value = 12900
$ws.Range("L132").Value = 19322.637
$ws.Range("M132").Value = -10370
$ws.Range("N132").Value = -24382.637
$ws.Range("H140").Value = 68533.45
$ws.Range("J140").Value = 68533.45
$ws.Range("L140").Value = 68533.45
$ws.Range("N140").Value = -78893.45

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8992.6
$ws.Range("I136").Value = 4793.5
$ws.Range("J136").Value = 9914.353999999999
$ws.Range("K136").Value = 14380.5
$ws.Range("L136").Value = 29743.062
$ws.Range("M136").Value = -11830.5
$ws.Range("N136").Value = -34843.062
